# Updates the cryptos list: Price (col D) / Volume(1h) (col E) changes, and the
# swap of the Maker / Bittensor rows (44-45) with their own updated values.
#
# Cells whose new text would otherwise be auto-recognized by Excel as a number
# (e.g. "580.16") are entered with a leading apostrophe so they are stored as text
# (General format) just like the rest of the column, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.810.19"
$ws.Range("E2").Value = "  +4.74%  "
$ws.Range("D3").Value = "3.263.23"
$ws.Range("E3").Value = "  +4.82%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'580.16"
$ws.Range("E5").Value = "  +2.72%  "
$ws.Range("D6").Value = "'181.85"
$ws.Range("E6").Value = "  +8.62%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").Value = "3.263.75"
$ws.Range("E9").Value = "  +4.87%  "
$ws.Range("E10").Value = "  +9.68%  "
$ws.Range("E11").Value = "  +3.62%  "
$ws.Range("D12").Value = "'0.417"
$ws.Range("E12").Value = "  +8.40%  "
$ws.Range("D13").Value = "3.833.79"
$ws.Range("E13").Value = "  +5.10%  "
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").Value = "'28.45"
$ws.Range("E15").Value = "  +6.90%  "
$ws.Range("D16").Value = "67.771.22"
$ws.Range("E16").Value = "  +4.72%  "
$ws.Range("E17").Value = "  +5.34%  "
$ws.Range("D18").Value = "3.277.19"
$ws.Range("E18").Value = "  +5.26%  "
$ws.Range("E19").Value = "  +4.46%  "
$ws.Range("E20").Value = "  +7.57%  "
$ws.Range("D21").Value = "'375.49"
$ws.Range("E21").Value = "  +6.93%  "
$ws.Range("E22").Value = "  +7.43%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'71.00"
$ws.Range("E24").Value = "  +4.21%  "
$ws.Range("E25").Value = "  +5.07%  "
$ws.Range("E26").Value = "  +8.69%  "
$ws.Range("D27").Value = "'9.61"
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("E28").Value = "  +4.08%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  +4.83%  "
$ws.Range("E31").Value = "  +9.62%  "
$ws.Range("D32").Value = "'22.75"
$ws.Range("E32").Value = "  +5.76%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +9.22%  "
$ws.Range("D35").Value = "'6.93"
$ws.Range("E35").Value = "  +6.81%  "
$ws.Range("D36").Value = "'163.72"
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("E37").Value = "  +7.21%  "
$ws.Range("D38").Value = "'0.849"
$ws.Range("E38").Value = "  +4.35%  "
$ws.Range("E39").Value = "  +7.25%  "
$ws.Range("D40").Value = "'6.82"
$ws.Range("E40").Value = "  +14.16%  "
$ws.Range("E41").Value = "  +4.05%  "
$ws.Range("D42").Value = "'4.66"
$ws.Range("E42").Value = "  +13.98%  "
$ws.Range("E43").Value = "  +9.62%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "'354.83"
$ws.Range("E44").Value = "  +13.05%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.703.10"
$ws.Range("E45").Value = "  +3.36%  "
$ws.Range("D46").Value = "'25.40"
$ws.Range("E46").Value = "  +8.69%  "
$ws.Range("D47").Value = "'40.86"
$ws.Range("E47").Value = "  +4.13%  "
$ws.Range("D48").Value = "'0.0681"
$ws.Range("E48").Value = "  +5.70%  "
$ws.Range("E49").Value = "  +5.10%  "
$ws.Range("E50").Value = "  +8.39%  "
$ws.Range("E51").Value = "  +1.52%  "
